# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) is recomputed from the newly-regenerated s_vals and
# rewritten for every data row. Row 37 also picks up a corrected IP/IF
# (columns H/J) recompute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated K values for rows 2..49 (row 19 unchanged / already 0)
$kVals = @{
    2  = 3
    3  = 3
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 3
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 4
    14 = 3
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 1
    25 = 4
    26 = 2
    27 = 5
    28 = 1
    29 = 0
    30 = 4
    31 = 2
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 2
    37 = 4
    38 = 0
    39 = 6
    40 = 1
    41 = 0
    42 = 1
    43 = 2
    44 = 2
    45 = 3
    46 = 2
    47 = 6
    48 = 2
    49 = 0
}

foreach ($row in $kVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $kVals[$row]
}

# Row 37 also regenerates IP (H) / IF (J) alongside its K value
$ws.Cells.Item(37, 8).Value = 6
$ws.Cells.Item(37, 10).Value = 8
